$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iteration 5")

# ---------------------------------------------------------------------------
# 1. New bug rows (9-13), mirroring the layout of row 8 (S/N, Functionality,
#    Files Affected, Description of bug, Date Found, Points, Severity,
#    Status columns B..I). Columns J,K,L stay blank but keep formatted style.
# ---------------------------------------------------------------------------

$rows = @(
    @{ Row = 9;  SN = 2; Func = "Temperature Charts";     Files = "viewPatientInformation.jsp"; Desc = "Date and time stated is wrong" },
    @{ Row = 10; SN = 3; Func = "Respiratory Rate Chart";  Files = "viewPatientInformation.jsp"; Desc = "Date and time stated is wrong" },
    @{ Row = 11; SN = 4; Func = "Heart Rate Chart";        Files = "viewPatientInformation.jsp"; Desc = "Date and time stated is wrong" },
    @{ Row = 12; SN = 5; Func = "Blood Pressure Chart";    Files = "viewPatientInformation.jsp"; Desc = "Date and time stated is wrong" },
    @{ Row = 13; SN = 6; Func = "SPO Chart";               Files = "viewPatientInformation.jsp"; Desc = "Date and time stated is wrong" }
)

foreach ($item in $rows) {
    $r = $item.Row

    $ws.Range("B${r}").Value = $item.SN
    $ws.Range("C${r}").Value = $item.Func
    $ws.Range("D${r}").Value = $item.Files
    $ws.Range("E${r}").Value = $item.Desc
    $ws.Range("F${r}").Value = 42001
    $ws.Range("G${r}").Value = 5
    $ws.Range("H${r}").Value = "High Impact "
    $ws.Range("I${r}").Value = "Unsolved"
    $ws.Range("J${r}").Value = ""
    $ws.Range("K${r}").Value = ""
    $ws.Range("L${r}").Value = ""

    # Copy the existing formatting already present in row 8 so no stray new
    # styles get invented; only B8 (below) needs an actually-new style.
    $ws.Range("C8").Copy()
    $ws.Range("B${r}").PasteSpecial(-4122)

    $ws.Range("D8").Copy()
    $ws.Range("C${r}").PasteSpecial(-4122)

    $ws.Range("E8").Copy()
    $ws.Range("D${r}").PasteSpecial(-4122)

    $ws.Range("F8").Copy()
    $ws.Range("E${r}:F${r}").PasteSpecial(-4122)

    $ws.Range("G8").Copy()
    $ws.Range("G${r}").PasteSpecial(-4122)

    $ws.Range("I8").Copy()
    $ws.Range("H${r}:I${r}").PasteSpecial(-4122)

    $ws.Range("K8").Copy()
    $ws.Range("J${r}").PasteSpecial(-4122)

    $ws.Range("L8").Copy()
    $ws.Range("K${r}:L${r}").PasteSpecial(-4122)
}

# Row heights: 9 & 10 keep their original (taller) custom heights, 11-13 get
# the new shorter height that Excel applied when the table grew.
$ws.Rows(9).RowHeight = 54.75
$ws.Rows(10).RowHeight = 57
$ws.Rows(11).RowHeight = 27.75
$ws.Rows(12).RowHeight = 27.75
$ws.Rows(13).RowHeight = 27.75

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. B8: give the S/N cell the new centered (no-wrap) style.
# ---------------------------------------------------------------------------
$b8 = $ws.Range("B8")
$b8.HorizontalAlignment = -4108
$b8.VerticalAlignment = -4108
$b8.WrapText = $false
$b8.Font.Name = "Century Gothic"
$b8.Font.Size = 10
$b8.Interior.Pattern = 1
$b8.Interior.ThemeColor = 0

# ---------------------------------------------------------------------------
# 3. D4 total formula now spans the whole bug list.
# ---------------------------------------------------------------------------
$ws.Range("D4").Formula = "=SUM(G8:G116)"

# ---------------------------------------------------------------------------
# 4. Selection moves to H11.
# ---------------------------------------------------------------------------
$ws.Range("H11").Select()

$wb.Application.CalculateFull()
